# Update TestLoginData.xlsx: refresh the jsessionid URLs recorded for the
# two "Đăng nhập" rows (rows 3 and 4 -> cells E3 and E4 on Sheet1) with new
# session URLs captured from a re-run of the automation tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = "http://localhost:8080/home/index;jsessionid=FF060E78F79C609A0474E1C28AD3349E"
$ws.Range("E4").Value = "http://localhost:8080/home/index;jsessionid=5B0E62695885E2E489F4764B6A4AD689"
